$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to update: maps cell address -> new value.
# Values that look like plain numbers (and would otherwise be
# reinterpreted by Excel as numeric, losing formatting such as
# trailing zeros) are written with a forced Text number format
# and then have their style reset back to Normal so no stray
# explicit style index is left behind on the cell.
$updates = [ordered]@{
    "D2" = "56.825.43"
    "E2" = "  +4.02%  "
    "D3" = "3.247.58"
    "E3" = "  +1.83%  "
    "E4" = "  +0.09%  "
    "D5" = "396.06"
    "E5" = "  -1.75%  "
    "D6" = "108.60"
    "E6" = "  -0.41%  "
    "D7" = "0.582"
    "E7" = "  +5.47%  "
    "D8" = "3.241.52"
    "E8" = "  +1.88%  "
    "D9" = "1.00"
    "E9" = "  +0.05%  "
    "D10" = "0.624"
    "E10" = "  +0.40%  "
    "D11" = "39.13"
    "E11" = "  +0.33%  "
    "D12" = "0.0987"
    "E12" = "  +11.86%  "
    "E13" = "  +1.75%  "
    "D14" = "3.761.58"
    "E14" = "  +2.20%  "
    "D15" = "8.32"
    "E15" = "  +3.31%  "
    "D16" = "19.10"
    "E16" = "  -0.36%  "
    "D17" = "3.246.82"
    "E17" = "  +1.59%  "
    "E18" = "  -3.10%  "
    "D19" = "10.73"
    "E19" = "  +2.30%  "
    "D20" = "56.834.24"
    "E20" = "  +4.28%  "
    "D21" = "3.35"
    "E21" = "  +1.22%  "
    "D23" = "12.92"
    "E23" = "  +0.22%  "
    "D24" = "295.59"
    "E24" = "  +7.69%  "
    "D25" = "74.22"
    "E25" = "  +3.30%  "
    "D26" = "3.17"
    "E26" = "  -3.46%  "
    "D27" = "28.13"
    "E27" = "  +1.32%  "
    "D28" = "4.35"
    "E28" = "  +1.03%  "
    "D29" = "7.61"
    "E29" = "  -4.85%  "
    "D30" = "7.25"
    "E30" = "  -1.81%  "
    "E31" = "  -1.53%  "
    "E32" = "  -0.01%  "
    "D33" = "11.27"
    "E33" = "  +1.48%  "
    "E34" = "  -4.17%  "
    "D35" = "39.39"
    "E35" = "  +6.94%  "
    "D36" = "0.0483"
    "E36" = "  -3.42%  "
    "E37" = "  +2.14%  "
    "D38" = "51.40"
    "E38" = "  +1.35%  "
    "D39" = "0.998"
    "E39" = "  -0.09%  "
    "E40" = "  -5.62%  "
    "D41" = "2.90"
    "E41" = "  +2.13%  "
    "D42" = "135.49"
    "E42" = "  +3.87%  "
    "E43" = "  +3.97%  "
    "D44" = "1.89"
    "E44" = "  -2.82%  "
    "B45" = "Celestia"
    "C45" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D45" = "17.05"
    "E45" = "  -1.38%  "
    "B46" = "NEARProtocol"
    "C46" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D46" = "3.95"
    "E46" = "  -4.22%  "
    "D47" = "0.279"
    "E47" = "  -4.12%  "
    "D48" = "22.15"
    "E48" = "  -0.59%  "
    "E49" = "  +3.22%  "
    "D50" = "2.151.50"
    "E50" = "  +2.95%  "
    "D51" = "2.33"
    "E51" = "  -7.21%  "
}

$numericLike = @(
    "D5"
    "D6"
    "D7"
    "D9"
    "D10"
    "D11"
    "D12"
    "D15"
    "D16"
    "D19"
    "D21"
    "D23"
    "D24"
    "D25"
    "D26"
    "D27"
    "D28"
    "D29"
    "D30"
    "D33"
    "D35"
    "D36"
    "D38"
    "D39"
    "D41"
    "D42"
    "D44"
    "D45"
    "D46"
    "D47"
    "D48"
    "D51"
)

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    if ($numericLike -contains $addr) {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$addr]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$addr]
    }
}
